# Update cryptos list values (prices & 1h volume changes) to the latest snapshot.
# Rows 38-41 and 49-51: coin name/link/price/volume cells re-ordered (ranking shuffled).
# Note: Price (column D) values are stored as plain text in the source sheet even when
# they look numeric (e.g. "9.40", "0.0000282"), so a leading apostrophe is used to force
# Excel to keep them as text instead of auto-converting to numbers, matching the original
# inlineStr formatting without altering the cell's number format/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


$ws.Range("D2").Value = '67.157.47'
$ws.Range("E2").Value = '  +3.68%  '
$ws.Range("D3").Value = '3.451.54'
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '''570.44'
$ws.Range("E5").Value = '  +2.18%  '
$ws.Range("D6").Value = '''185.67'
$ws.Range("E6").Value = '  +6.24%  '
$ws.Range("E7").Value = '  +1.86%  '
$ws.Range("D8").Value = '3.445.21'
$ws.Range("E8").Value = '  +2.32%  '
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("D10").Value = '''0.176'
$ws.Range("E10").Value = '  +6.60%  '
$ws.Range("D11").Value = '''0.645'
$ws.Range("E11").Value = '  +2.35%  '
$ws.Range("D12").Value = '''55.47'
$ws.Range("E12").Value = '  +2.77%  '
$ws.Range("D13").Value = '''0.0000282'
$ws.Range("E13").Value = '  +2.49%  '
$ws.Range("D14").Value = '''9.40'
$ws.Range("E14").Value = '  +3.49%  '
$ws.Range("D15").Value = '3.988.43'
$ws.Range("E15").Value = '  +1.71%  '
$ws.Range("D16").Value = '''18.61'
$ws.Range("E16").Value = '  +1.98%  '
$ws.Range("D17").Value = '3.436.88'
$ws.Range("E17").Value = '  +1.59%  '
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").Value = '66.921.64'
$ws.Range("E19").Value = '  +3.63%  '
$ws.Range("D20").Value = '''12.04'
$ws.Range("E20").Value = '  +2.49%  '
$ws.Range("E21").Value = '  +2.19%  '
$ws.Range("D22").Value = '''479.83'
$ws.Range("E22").Value = '  +4.11%  '
$ws.Range("D23").Value = '''4.98'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").Value = '''15.25'
$ws.Range("E24").Value = '  +13.23%  '
$ws.Range("D25").Value = '''4.19'
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("D26").Value = '''89.74'
$ws.Range("E26").Value = '  +3.97%  '
$ws.Range("D27").Value = '''2.97'
$ws.Range("E27").Value = '  +1.10%  '
$ws.Range("D28").Value = '''11.08'
$ws.Range("E28").Value = '  +2.72%  '
$ws.Range("D29").Value = '''8.94'
$ws.Range("E29").Value = '  +2.51%  '
$ws.Range("D30").Value = '''31.62'
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("D31").Value = '''6.97'
$ws.Range("E31").Value = '  +4.04%  '
$ws.Range("D32").Value = '''11.66'
$ws.Range("E32").Value = '  +1.85%  '
$ws.Range("D33").Value = '''590.48'
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("D34").Value = '''63.19'
$ws.Range("E34").Value = '  +3.37%  '
$ws.Range("E35").Value = '  +1.99%  '
$ws.Range("E36").Value = '  +6.74%  '
$ws.Range("B38").Value = 'PEPE'
$ws.Range("C38").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D38").Value = '0.0₃0798'
$ws.Range("E38").Value = '  +7.79%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''3.63'
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = '''0.392'
$ws.Range("E40").Value = '  +6.35%  '
$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").Value = '''36.67'
$ws.Range("E41").Value = '  +3.56%  '
$ws.Range("D42").Value = '3.128.29'
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("D43").Value = '''2.93'
$ws.Range("E43").Value = '  +3.45%  '
$ws.Range("D44").Value = '''2.65'
$ws.Range("E44").Value = '  +8.58%  '
$ws.Range("D45").Value = '''0.0426'
$ws.Range("E45").Value = '  +2.91%  '
$ws.Range("D46").Value = '''2.81'
$ws.Range("E46").Value = '  +21.60%  '
$ws.Range("D47").Value = '''3.25'
$ws.Range("E47").Value = '  +4.21%  '
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '''8.71'
$ws.Range("E49").Value = '  +6.27%  '
$ws.Range("D50").Value = '''142.26'
$ws.Range("E50").Value = '  +2.04%  '
$ws.Range("B51").Value = 'FirstDigitalUSD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D51").Value = '''0.998'
$ws.Range("E51").Value = '  -0.14%  '
